# Auto-generated edit script applying the scheduled-runner value updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 699664.3
$ws.Range("J17").Value = 724464.8
$ws.Range("L17").Value = 2173394.4
$ws.Range("N17").Value = -2173730.4
$ws.Range("H32").Value = 800
$ws.Range("J32").Value = 800
$ws.Range("L32").Value = 800
$ws.Range("N32").Value = -1452
$ws.Range("H33").Value = 1341.75
$ws.Range("I33").Value = 956
$ws.Range("K33").Value = 956
$ws.Range("M33").Value = -727
$ws.Range("H64").Value = 6727
$ws.Range("I64").Value = 4499.5
$ws.Range("J64").Value = 7999.857
$ws.Range("K64").Value = 4499.5
$ws.Range("L64").Value = 7999.857
$ws.Range("M64").Value = -4251.5
$ws.Range("N64").Value = -8495.857
$ws.Range("H67").Value = 6727
$ws.Range("I67").Value = 4499.5
$ws.Range("J67").Value = 7999.857
$ws.Range("K67").Value = 4499.5
$ws.Range("L67").Value = 7999.857
$ws.Range("M67").Value = -3641.5
$ws.Range("N67").Value = -9715.857
$ws.Range("H69").Value = 7182.4546
$ws.Range("I69").Value = 6997
$ws.Range("J69").Value = 7201
$ws.Range("K69").Value = 20991
$ws.Range("L69").Value = 21603
$ws.Range("M69").Value = -20117
$ws.Range("N69").Value = -23351
$ws.Range("H72").Value = 7182.4546
$ws.Range("I72").Value = 6997
$ws.Range("J72").Value = 7201
$ws.Range("K72").Value = 62973
$ws.Range("L72").Value = 64809
$ws.Range("M72").Value = -58605
$ws.Range("N72").Value = -73545
$ws.Range("H74").Value = 10215.25
$ws.Range("I74").Value = 10430.5
$ws.Range("K74").Value = 10430.5
$ws.Range("M74").Value = -9494.5
$ws.Range("H76").Value = 4499
$ws.Range("J76").Value = 4499
$ws.Range("L76").Value = 4499
$ws.Range("N76").Value = -5129
$ws.Range("H77").Value = 10215.25
$ws.Range("I77").Value = 10430.5
$ws.Range("K77").Value = 52152.5
$ws.Range("M77").Value = -47472.5
$ws.Range("H79").Value = 4499
$ws.Range("J79").Value = 4499
$ws.Range("L79").Value = 4499
$ws.Range("N79").Value = -6683
$ws.Range("H86").Value = 69515780
$ws.Range("I86").Value = 19153.154
$ws.Range("K86").Value = 19153.154
$ws.Range("M86").Value = -18030.154
$ws.Range("H89").Value = 69515780
$ws.Range("I89").Value = 19153.154
$ws.Range("K89").Value = 95765.76999999999
$ws.Range("M89").Value = -90149.76999999999
$ws.Range("H98").Value = 2575.7693
$ws.Range("I98").Value = 1698.375
$ws.Range("K98").Value = 1698.375
$ws.Range("M98").Value = -200.375
$ws.Range("H112").Value = 1719602.2
$ws.Range("I112").Value = 3052.625
$ws.Range("J112").Value = 2863968.5
$ws.Range("K112").Value = 9157.875
$ws.Range("L112").Value = 8591905.5
$ws.Range("M112").Value = -8049.875
$ws.Range("N112").Value = -8594121.5
$ws.Range("H113").Value = 8629.375
$ws.Range("J113").Value = 3308.2
$ws.Range("L113").Value = 3308.2
$ws.Range("N113").Value = -9816.200000000001
$ws.Range("H116").Value = 35732284
$ws.Range("I116").Value = 125045496
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 125045496
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -125042054
$ws.Range("N116").Value = -13884
$ws.Range("H122").Value = 2575.7693
$ws.Range("I122").Value = 1698.375
$ws.Range("K122").Value = 5095.125
$ws.Range("M122").Value = -2645.125
$ws.Range("H132").Value = 1140.7878
$ws.Range("I132").Value = 1168.1451
$ws.Range("K132").Value = 3504.4353
$ws.Range("M132").Value = -974.4353000000001
$ws.Range("H135").Value = 1173.909
$ws.Range("J135").Value = 1700
$ws.Range("L135").Value = 15300
$ws.Range("N135").Value = -20370
$ws.Range("H138").Value = 3656.8652
$ws.Range("I138").Value = 1480.3529
$ws.Range("J138").Value = 4170.7637
$ws.Range("K138").Value = 4441.0587
$ws.Range("L138").Value = 12512.2911
$ws.Range("M138").Value = 698.9412999999995
$ws.Range("N138").Value = -22792.2911
$ws.Range("H141").Value = 1025
$ws.Range("I141").Value = 1025
$ws.Range("K141").Value = 3075
$ws.Range("M141").Value = 2105

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17546092
$ws.Range("I32").Value = 19232494
$ws.Range("K32").Value = 19232494
$ws.Range("M32").Value = -19232207
$ws.Range("H63").Value = 4797.5
$ws.Range("I63").Value = 4371.25
$ws.Range("J63").Value = 5650
$ws.Range("K63").Value = 4371.25
$ws.Range("L63").Value = 5650
$ws.Range("M63").Value = -3685.25
$ws.Range("N63").Value = -7022
$ws.Range("H66").Value = 4797.5
$ws.Range("I66").Value = 4371.25
$ws.Range("J66").Value = 5650
$ws.Range("K66").Value = 21856.25
$ws.Range("L66").Value = 28250
$ws.Range("M66").Value = -18424.25
$ws.Range("N66").Value = -35114
$ws.Range("H74").Value = 2680.932
$ws.Range("I74").Value = 2210.8108
$ws.Range("J74").Value = 5165.857
$ws.Range("K74").Value = 2210.8108
$ws.Range("L74").Value = 5165.857
$ws.Range("M74").Value = -1336.8108
$ws.Range("N74").Value = -6913.857
$ws.Range("H77").Value = 2680.932
$ws.Range("I77").Value = 2210.8108
$ws.Range("J77").Value = 5165.857
$ws.Range("K77").Value = 11054.054
$ws.Range("L77").Value = 25829.285
$ws.Range("M77").Value = -6686.054
$ws.Range("N77").Value = -34565.285
$ws.Range("H81").Value = 120982
$ws.Range("J81").Value = 120982
$ws.Range("L81").Value = 120982
$ws.Range("N81").Value = -122978
$ws.Range("H84").Value = 120982
$ws.Range("J84").Value = 120982
$ws.Range("L84").Value = 362946
$ws.Range("N84").Value = -372930
$ws.Range("H118").Value = 118732.336
$ws.Range("J118").Value = 118732.336
$ws.Range("L118").Value = 118732.336
$ws.Range("N118").Value = -122046.336
$ws.Range("H122").Value = 5568
$ws.Range("I122").Value = 3946.6667
$ws.Range("K122").Value = 11840.0001
$ws.Range("M122").Value = -9390.000100000001
$ws.Range("H128").Value = 36180.4
$ws.Range("J128").Value = 36180.4
$ws.Range("L128").Value = 36180.4
$ws.Range("N128").Value = -46140.4
$ws.Range("H132").Value = 2721.5366
$ws.Range("I132").Value = 2355.0293
$ws.Range("K132").Value = 7065.0879
$ws.Range("M132").Value = -4535.0879
$ws.Range("H139").Value = 76299.60000000001
$ws.Range("J139").Value = 77624.5
$ws.Range("L139").Value = 77624.5
$ws.Range("N139").Value = -87904.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4124.8335
$ws.Range("I86").Value = 3949.8
$ws.Range("K86").Value = 3949.8
$ws.Range("M86").Value = -2826.8
$ws.Range("H89").Value = 4124.8335
$ws.Range("I89").Value = 3949.8
$ws.Range("K89").Value = 19749
$ws.Range("M89").Value = -14133
$ws.Range("H94").Value = 1209.9131
$ws.Range("I94").Value = 1024.0834
$ws.Range("K94").Value = 1024.0834
$ws.Range("M94").Value = -573.0834
$ws.Range("H107").Value = 2298.5789
$ws.Range("I107").Value = 1617.909
$ws.Range("K107").Value = 1617.909
$ws.Range("M107").Value = 302.0909999999999
$ws.Range("H116").Value = 99871
$ws.Range("J116").Value = 99871
$ws.Range("L116").Value = 99871
$ws.Range("N116").Value = -109049
$ws.Range("H134").Value = 4793.5835
$ws.Range("I134").Value = 2863.625
$ws.Range("K134").Value = 8590.875
$ws.Range("M134").Value = -6055.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 739.4545000000001
$ws.Range("J22").Value = 614.3333
$ws.Range("L22").Value = 614.3333
$ws.Range("N22").Value = -1314.3333
$ws.Range("H31").Value = 5519.3687
$ws.Range("I31").Value = 1924.4667
$ws.Range("J31").Value = 7863.8696
$ws.Range("K31").Value = 1924.4667
$ws.Range("L31").Value = 7863.8696
$ws.Range("M31").Value = -1629.4667
$ws.Range("N31").Value = -8453.8696
$ws.Range("H34").Value = 5519.3687
$ws.Range("I34").Value = 1924.4667
$ws.Range("J34").Value = 7863.8696
$ws.Range("K34").Value = 1924.4667
$ws.Range("L34").Value = 7863.8696
$ws.Range("M34").Value = -1722.4667
$ws.Range("N34").Value = -8267.8696
$ws.Range("H39").Value = 2508.3333
$ws.Range("I39").Value = 5525
$ws.Range("K39").Value = 5525
$ws.Range("M39").Value = -5134
$ws.Range("H49").Value = 2508.3333
$ws.Range("I49").Value = 5525
$ws.Range("K49").Value = 5525
$ws.Range("M49").Value = -5343
$ws.Range("H58").Value = 3588.524
$ws.Range("I58").Value = 3567.95
$ws.Range("K58").Value = 3567.95
$ws.Range("M58").Value = -3364.95
$ws.Range("H62").Value = 3995
$ws.Range("I62").Value = 3995
$ws.Range("K62").Value = 3995
$ws.Range("M62").Value = -3371
$ws.Range("H64").Value = 38000
$ws.Range("J64").Value = 38000
$ws.Range("L64").Value = 38000
$ws.Range("N64").Value = -38496
$ws.Range("H65").Value = 3995
$ws.Range("I65").Value = 3995
$ws.Range("K65").Value = 19975
$ws.Range("M65").Value = -16855
$ws.Range("H67").Value = 38000
$ws.Range("J67").Value = 38000
$ws.Range("L67").Value = 38000
$ws.Range("N67").Value = -39716
$ws.Range("H132").Value = 3377.4285
$ws.Range("I132").Value = 3515.8333
$ws.Range("K132").Value = 10547.4999
$ws.Range("M132").Value = -8017.499899999999
$ws.Range("H134").Value = 1969.75
$ws.Range("I134").Value = 1658.3572
$ws.Range("K134").Value = 4975.071599999999
$ws.Range("M134").Value = -2440.071599999999
$ws.Range("H136").Value = 3588.524
$ws.Range("I136").Value = 3567.95
$ws.Range("K136").Value = 10703.85
$ws.Range("M136").Value = -8153.849999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4980
$ws.Range("I3").Value = 5850
$ws.Range("K3").Value = 17550
$ws.Range("M3").Value = -17438
$ws.Range("H26").Value = 66.666664
$ws.Range("J26").Value = 50
$ws.Range("L26").Value = 150
$ws.Range("N26").Value = -726
$ws.Range("H32").Value = 111111550
$ws.Range("J32").Value = 500000260
$ws.Range("L32").Value = 1500000780
$ws.Range("N32").Value = -1500001346
$ws.Range("H39").Value = 4991.4
$ws.Range("J39").Value = 4991.4
$ws.Range("L39").Value = 14974.2
$ws.Range("N39").Value = -15562.2
$ws.Range("H47").Value = 722.8
$ws.Range("I47").Value = 789.6667
$ws.Range("J47").Value = 622.5
$ws.Range("K47").Value = 2369.0001
$ws.Range("L47").Value = 1867.5
$ws.Range("M47").Value = -1938.0001
$ws.Range("N47").Value = -2729.5
$ws.Range("H107").Value = 787.3333
$ws.Range("J107").Value = 619.625
$ws.Range("L107").Value = 1858.875
$ws.Range("N107").Value = -5698.875
$ws.Range("H113").Value = 1044.5834
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1044.5834
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3133.7502
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -7473.7502
$ws.Range("H121").Value = 4774640.5
$ws.Range("I121").Value = 761.8333
$ws.Range("J121").Value = 6977969
$ws.Range("K121").Value = 2285.4999
$ws.Range("L121").Value = 20933907
$ws.Range("M121").Value = -975.4998999999998
$ws.Range("N121").Value = -20936527
$ws.Range("H129").Value = 2022.4445
$ws.Range("I129").Value = 784.5
$ws.Range("K129").Value = 2353.5
$ws.Range("M129").Value = 2646.5
$ws.Range("H134").Value = 5301.5713
$ws.Range("I134").Value = 5301.5713
$ws.Range("K134").Value = 15904.7139
$ws.Range("M134").Value = -10834.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2449.75
$ws.Range("I80").Value = 1933
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 1933
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -935
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 2449.75
$ws.Range("I83").Value = 1933
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 9665
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -4673
$ws.Range("N83").Value = -29984
$ws.Range("H103").Value = 61499.5
$ws.Range("J103").Value = 61499.5
$ws.Range("L103").Value = 61499.5
$ws.Range("N103").Value = -63843.5
$ws.Range("H107").Value = 1457.5
$ws.Range("I107").Value = 1171.2
$ws.Range("J107").Value = 1662
$ws.Range("K107").Value = 1171.2
$ws.Range("L107").Value = 1662
$ws.Range("M107").Value = 748.8
$ws.Range("N107").Value = -5502
$ws.Range("H122").Value = 3429.5
$ws.Range("I122").Value = 3429.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10288.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -7838.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3592.4285
$ws.Range("I22").Value = 3279.7
$ws.Range("J22").Value = 3876.7273
$ws.Range("K22").Value = 3279.7
$ws.Range("L22").Value = 3876.7273
$ws.Range("M22").Value = -2984.7
$ws.Range("N22").Value = -4466.7273
$ws.Range("H27").Value = 3592.4285
$ws.Range("I27").Value = 3279.7
$ws.Range("J27").Value = 3876.7273
$ws.Range("K27").Value = 3279.7
$ws.Range("L27").Value = 3876.7273
$ws.Range("M27").Value = -3172.7
$ws.Range("N27").Value = -4090.7273
$ws.Range("H38").Value = 45032.5
$ws.Range("I38").Value = 30030
$ws.Range("J38").Value = 60035
$ws.Range("K38").Value = 30030
$ws.Range("L38").Value = 60035
$ws.Range("M38").Value = -29620
$ws.Range("N38").Value = -60855
$ws.Range("H40").Value = 33340404
$ws.Range("I40").Value = 66671468
$ws.Range("K40").Value = 66671468
$ws.Range("M40").Value = -66671332
$ws.Range("H82").Value = 18762.334
$ws.Range("I82").Value = 2479.4285
$ws.Range("J82").Value = 41558.4
$ws.Range("K82").Value = 2479.4285
$ws.Range("L82").Value = 41558.4
$ws.Range("M82").Value = -2118.4285
$ws.Range("N82").Value = -42280.4
$ws.Range("H85").Value = 18762.334
$ws.Range("I85").Value = 2479.4285
$ws.Range("J85").Value = 41558.4
$ws.Range("K85").Value = 2479.4285
$ws.Range("L85").Value = 41558.4
$ws.Range("M85").Value = -1231.4285
$ws.Range("N85").Value = -44054.4
$ws.Range("H93").Value = 100002290
$ws.Range("I93").Value = 166667890
$ws.Range("J93").Value = 3875
$ws.Range("K93").Value = 166667890
$ws.Range("L93").Value = 3875
$ws.Range("M93").Value = -166666642
$ws.Range("N93").Value = -6371
$ws.Range("H107").Value = 9997
$ws.Range("I107").Value = 9997
$ws.Range("K107").Value = 9997
$ws.Range("M107").Value = -8077
$ws.Range("H122").Value = 14591.667
$ws.Range("I122").Value = 6888
$ws.Range("K122").Value = 20664
$ws.Range("M122").Value = -18214
$ws.Range("H132").Value = 5571.838
$ws.Range("I132").Value = 5609.9
$ws.Range("J132").Value = 5408.7144
$ws.Range("K132").Value = 16829.7
$ws.Range("L132").Value = 16226.1432
$ws.Range("M132").Value = -14299.7
$ws.Range("N132").Value = -21286.1432
$ws.Range("H136").Value = 2078.7778
$ws.Range("I136").Value = 1451.1666
$ws.Range("K136").Value = 4353.4998
$ws.Range("M136").Value = -1803.4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 8253.4
$ws.Range("I32").Value = 8253.4
$ws.Range("K32").Value = 8253.4
$ws.Range("M32").Value = -7936.4
$ws.Range("H75").Value = 79952.5
$ws.Range("I75").Value = 79960
$ws.Range("J75").Value = 79945
$ws.Range("K75").Value = 79960
$ws.Range("L75").Value = 79945
$ws.Range("M75").Value = -79024
$ws.Range("N75").Value = -81817
$ws.Range("H78").Value = 79952.5
$ws.Range("I78").Value = 79960
$ws.Range("J78").Value = 79945
$ws.Range("K78").Value = 239880
$ws.Range("L78").Value = 239835
$ws.Range("M78").Value = -235200
$ws.Range("N78").Value = -249195
$ws.Range("H94").Value = 59666.668
$ws.Range("J94").Value = 59666.668
$ws.Range("L94").Value = 59666.668
$ws.Range("N94").Value = -61468.668
$ws.Range("H107").Value = 594.58826
$ws.Range("I107").Value = 467.45456
$ws.Range("J107").Value = 827.6667
$ws.Range("K107").Value = 1402.36368
$ws.Range("L107").Value = 2483.0001
$ws.Range("M107").Value = 517.6363200000001
$ws.Range("N107").Value = -6323.0001
$ws.Range("H122").Value = 250012350
$ws.Range("I122").Value = 333344900
$ws.Range("K122").Value = 1000034700
$ws.Range("M122").Value = -1000032250
$ws.Range("H132").Value = 2005.7959
$ws.Range("I132").Value = 2016.3334
$ws.Range("K132").Value = 6049.0002
$ws.Range("M132").Value = -3519.0002
$ws.Range("H136").Value = 44561.25
$ws.Range("I136").Value = 2167.2
$ws.Range("K136").Value = 6501.599999999999
$ws.Range("M136").Value = -3951.599999999999

